$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 4 new rows above row 72 (pushes the Class table from rows 72-77 down to 76-81)
$ws.Rows("72:75").Insert()

# Write the note into B71 (row just above the table)
$ws.Range("B71").Value = "Note: GW value assumes first usable address in range."

# Select B71 and scroll so A55 is the top-left visible cell (matches author's final view state)
$ws.Range("B71").Select()
$ws.Application.ActiveWindow.ScrollRow = 55
